$wb = $excel.ActiveWorkbook

$wsRun = $wb.Worksheets.Item("RunManager")
$wsScenarios = $wb.Worksheets.Item("TestcaseScenarios")
$wsData = $wb.Worksheets.Item("TestData")

# --- TestcaseScenarios sheet ---
# Column A header changes from "SerialNo" to "TestCaseName"
$wsScenarios.Range("A1").Value = "TestCaseName"

# Row 2 : TC_Login -> TC_Login1 (was numeric 1 before)
$wsScenarios.Range("A2").Value = "TC_Login1"

# New row 3 : TC_Login2
$wsScenarios.Range("A3").Value = "TC_Login2"

# --- RunManager sheet ---
# Row2 / Row3 TestCaseName column (D) now reference TC_Login1 / TC_Login2
$wsRun.Range("D2").Value = "TC_Login1"
$wsRun.Range("D3").Value = "TC_Login2"

# --- Selections / active sheet / views ---
$wsScenarios.Range("D20:E20").Select()
$wsRun.Range("D1:D1048576").Select()

$wsScenarios.Activate()

$wb.Save()
